{"js": "// Update the p-value table in Fig 3: replace the Cod/Hake trend-line\n// p-values while leaving every other value (and all formatting) intact.\n//\n// Table layout (row, col -- both 0-based via Table.getCell):\n//   row 0: \"\", \"Cod\", \"Hake\"                 (header)\n//   row 1: \"GDP 2016\",      0.76, 0.78\n//   row 2: \"OHI fisheries\", 0.29, 0.47\n//   row 3: \"OHI economic\",  0.95, 0.87\n//   row 4: \"Readiness\",     0.18, 0.45\n//   row 5: \"Vulnerability\", 0.02, 0.22  (Hake value unchanged)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// { row, col, old, new } pairs taken directly from the diff.\nconst updates = [\n  { row: 1, col: 1, from: \"0.76\", to: \"0.73\" }, // GDP 2016 / Cod\n  { row: 1, col: 2, from: \"0.78\", to: \"0.76\" }, // GDP 2016 / Hake\n  { row: 2, col: 1, from: \"0.29\", to: \"0.35\" }, // OHI fisheries / Cod\n  { row: 2, col: 2, from: \"0.47\", to: \"0.19\" }, // OHI fisheries / Hake\n  { row: 3, col: 1, from: \"0.95\", to: \"0.99\" }, // OHI economic / Cod\n  { row: 3, col: 2, from: \"0.87\", to: \"0.73\" }, // OHI economic / Hake\n  { row: 4, col: 1, from: \"0.18\", to: \"0.13\" }, // Readiness / Cod\n  { row: 4, col: 2, from: \"0.45\", to: \"0.34\" }, // Readiness / Hake\n  { row: 5, col: 1, from: \"0.02\", to: \"0.03\" }, // Vulnerability / Cod\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.load(\"value\");\n  await context.sync();\n\n  const range = cell.getRange();\n  if (cell.value === u.from) {\n    // Replace the whole cell's text run -- keeps the existing run\n    // formatting (font/size/color) on the surviving run.\n    range.insertText(u.to, \"Replace\");\n  } else {\n    // Fall back to a plain-text replace in case the cell already carries\n    // something other than the expected original value.\n    range.insertText(cell.value.split(u.from).join(u.to), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the p-value table in Fig 3: replace the Cod/Hake trend-line\n# p-values while leaving every other value (and all formatting) intact.\n#\n# Table layout (Cell(row, col) is 1-based, matching the Word object model):\n#   row 1: \"\", \"Cod\", \"Hake\"                 (header)\n#   row 2: \"GDP 2016\",      0.76, 0.78\n#   row 3: \"OHI fisheries\", 0.29, 0.47\n#   row 4: \"OHI economic\",  0.95, 0.87\n#   row 5: \"Readiness\",     0.18, 0.45\n#   row 6: \"Vulnerability\", 0.02, 0.22  (Hake value unchanged)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# { row, col, old, new } pairs taken directly from the diff.\n$updates = @(\n    @{ Row = 2; Col = 2; From = \"0.76\"; To = \"0.73\" },  # GDP 2016 / Cod\n    @{ Row = 2; Col = 3; From = \"0.78\"; To = \"0.76\" },  # GDP 2016 / Hake\n    @{ Row = 3; Col = 2; From = \"0.29\"; To = \"0.35\" },  # OHI fisheries / Cod\n    @{ Row = 3; Col = 3; From = \"0.47\"; To = \"0.19\" },  # OHI fisheries / Hake\n    @{ Row = 4; Col = 2; From = \"0.95\"; To = \"0.99\" },  # OHI economic / Cod\n    @{ Row = 4; Col = 3; From = \"0.87\"; To = \"0.73\" },  # OHI economic / Hake\n    @{ Row = 5; Col = 2; From = \"0.18\"; To = \"0.13\" },  # Readiness / Cod\n    @{ Row = 5; Col = 3; From = \"0.45\"; To = \"0.34\" },  # Readiness / Hake\n    @{ Row = 6; Col = 2; From = \"0.02\"; To = \"0.03\" }   # Vulnerability / Cod\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -eq $u.From) {\n        $cell.Range.Text = $u.To\n    } else {\n        $cell.Range.Text = $current.Replace($u.From, $u.To)\n    }\n}\n"}
